# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the existing header cells (bold, bordered,
# center/top aligned) by copying formatting from an existing header cell.
$headerSrc = $ws.Range("AC1")
$headerRange = $ws.Range("AD1:AF1")
$headerSrc.Copy()
$headerRange.PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins=63, Losses=97, Ties=1) for every data row.
$lastRow = 39
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 63
    $ws.Cells.Item($r, 31).Value = 97
    $ws.Cells.Item($r, 32).Value = 1
}

Write-Host "done"
